$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1006.25
$ws.Cells.Item(28, 9).Value = 446.41177
$ws.Cells.Item(28, 10).Value = 2365.8572
$ws.Cells.Item(28, 11).Value = 446.41177
$ws.Cells.Item(28, 12).Value = 2365.8572
$ws.Cells.Item(28, 13).Value = 38.58823000000001
$ws.Cells.Item(28, 14).Value = -3335.8572
$ws.Cells.Item(94, 8).Value = 3424.5557
$ws.Cells.Item(94, 9).Value = 2880.5
$ws.Cells.Item(94, 11).Value = 2880.5
$ws.Cells.Item(94, 13).Value = -2429.5
$ws.Cells.Item(97, 8).Value = 433.125
$ws.Cells.Item(97, 9).Value = 452.5
$ws.Cells.Item(97, 10).Value = 426.66666
$ws.Cells.Item(97, 11).Value = 1357.5
$ws.Cells.Item(97, 12).Value = 1279.99998
$ws.Cells.Item(97, 13).Value = -861.5
$ws.Cells.Item(97, 14).Value = -2271.99998
$ws.Cells.Item(99, 8).Value = 483.27777
$ws.Cells.Item(99, 10).Value = 726.375
$ws.Cells.Item(99, 12).Value = 2179.125
$ws.Cells.Item(99, 14).Value = -5175.125
$ws.Cells.Item(101, 8).Value = 4635.273
$ws.Cells.Item(101, 9).Value = 381
$ws.Cells.Item(101, 10).Value = 7770
$ws.Cells.Item(101, 11).Value = 1143
$ws.Cells.Item(101, 12).Value = 23310
$ws.Cells.Item(101, 13).Value = 479
$ws.Cells.Item(101, 14).Value = -26554
$ws.Cells.Item(132, 8).Value = 3361.0286
$ws.Cells.Item(132, 9).Value = 2843.7932
$ws.Cells.Item(132, 10).Value = 5861
$ws.Cells.Item(132, 11).Value = 8531.3796
$ws.Cells.Item(132, 12).Value = 17583
$ws.Cells.Item(132, 13).Value = -6001.3796
$ws.Cells.Item(132, 14).Value = -22643
$ws.Cells.Item(138, 8).Value = 1276.4445
$ws.Cells.Item(138, 9).Value = 1078.1818
$ws.Cells.Item(138, 11).Value = 3234.5454
$ws.Cells.Item(138, 13).Value = 1905.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12052860
$ws.Cells.Item(32, 9).Value = 14928112
$ws.Cells.Item(32, 10).Value = 12738.1875
$ws.Cells.Item(32, 11).Value = 14928112
$ws.Cells.Item(32, 12).Value = 12738.1875
$ws.Cells.Item(32, 13).Value = -14927825
$ws.Cells.Item(32, 14).Value = -13312.1875
$ws.Cells.Item(61, 8).Value = 1270.5636
$ws.Cells.Item(61, 9).Value = 901.2857
$ws.Cells.Item(61, 10).Value = 2463.6155
$ws.Cells.Item(61, 11).Value = 901.2857
$ws.Cells.Item(61, 12).Value = 2463.6155
$ws.Cells.Item(61, 13).Value = -689.2857
$ws.Cells.Item(61, 14).Value = -2887.6155
$ws.Cells.Item(97, 8).Value = 503.96155
$ws.Cells.Item(97, 9).Value = 392.16666
$ws.Cells.Item(97, 10).Value = 755.5
$ws.Cells.Item(97, 11).Value = 392.16666
$ws.Cells.Item(97, 12).Value = 755.5
$ws.Cells.Item(97, 13).Value = 103.83334
$ws.Cells.Item(97, 14).Value = -1747.5
$ws.Cells.Item(132, 8).Value = 1715.6285
$ws.Cells.Item(132, 9).Value = 1576.909
$ws.Cells.Item(132, 10).Value = 1950.3846
$ws.Cells.Item(132, 11).Value = 4730.727000000001
$ws.Cells.Item(132, 12).Value = 5851.1538
$ws.Cells.Item(132, 13).Value = -2200.727000000001
$ws.Cells.Item(132, 14).Value = -10911.1538
$ws.Cells.Item(136, 8).Value = 1270.5636
$ws.Cells.Item(136, 9).Value = 901.2857
$ws.Cells.Item(136, 10).Value = 2463.6155
$ws.Cells.Item(136, 11).Value = 2703.8571
$ws.Cells.Item(136, 12).Value = 7390.8465
$ws.Cells.Item(136, 13).Value = -153.8571000000002
$ws.Cells.Item(136, 14).Value = -12490.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2962.95
$ws.Cells.Item(20, 9).Value = 1681.6666
$ws.Cells.Item(20, 10).Value = 4884.875
$ws.Cells.Item(20, 11).Value = 1681.6666
$ws.Cells.Item(20, 12).Value = 4884.875
$ws.Cells.Item(20, 13).Value = -1434.6666
$ws.Cells.Item(20, 14).Value = -5378.875
$ws.Cells.Item(94, 8).Value = 405.58334
$ws.Cells.Item(94, 9).Value = 350.77777
$ws.Cells.Item(94, 10).Value = 570
$ws.Cells.Item(94, 11).Value = 350.77777
$ws.Cells.Item(94, 12).Value = 570
$ws.Cells.Item(94, 13).Value = 100.22223
$ws.Cells.Item(94, 14).Value = -1472
$ws.Cells.Item(99, 8).Value = 1548.1538
$ws.Cells.Item(99, 9).Value = 670
$ws.Cells.Item(99, 10).Value = 2745.6365
$ws.Cells.Item(99, 11).Value = 670
$ws.Cells.Item(99, 12).Value = 2745.6365
$ws.Cells.Item(99, 13).Value = 828
$ws.Cells.Item(99, 14).Value = -5741.636500000001
$ws.Cells.Item(107, 8).Value = 2664.3
$ws.Cells.Item(107, 9).Value = 2858.6667
$ws.Cells.Item(107, 10).Value = 2081.2
$ws.Cells.Item(107, 11).Value = 2858.6667
$ws.Cells.Item(107, 12).Value = 2081.2
$ws.Cells.Item(107, 13).Value = -938.6667000000002
$ws.Cells.Item(107, 14).Value = -5921.2
$ws.Cells.Item(134, 8).Value = 1728.0197
$ws.Cells.Item(134, 9).Value = 1476.262
$ws.Cells.Item(134, 10).Value = 2902.889
$ws.Cells.Item(134, 11).Value = 4428.786
$ws.Cells.Item(134, 12).Value = 8708.667000000001
$ws.Cells.Item(134, 13).Value = -1893.786
$ws.Cells.Item(134, 14).Value = -13778.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2024.7693
$ws.Cells.Item(16, 9).Value = 1101.8334
$ws.Cells.Item(16, 10).Value = 2815.8572
$ws.Cells.Item(16, 11).Value = 1101.8334
$ws.Cells.Item(16, 12).Value = 2815.8572
$ws.Cells.Item(16, 13).Value = -814.8334
$ws.Cells.Item(16, 14).Value = -3389.8572
$ws.Cells.Item(113, 8).Value = 2024.7693
$ws.Cells.Item(113, 9).Value = 1101.8334
$ws.Cells.Item(113, 10).Value = 2815.8572
$ws.Cells.Item(113, 11).Value = 1101.8334
$ws.Cells.Item(113, 12).Value = 2815.8572
$ws.Cells.Item(113, 13).Value = 1068.1666
$ws.Cells.Item(113, 14).Value = -7155.8572
$ws.Cells.Item(132, 8).Value = 1320.9574
$ws.Cells.Item(132, 9).Value = 1078.6571
$ws.Cells.Item(132, 10).Value = 2027.6666
$ws.Cells.Item(132, 11).Value = 3235.9713
$ws.Cells.Item(132, 12).Value = 6082.9998
$ws.Cells.Item(132, 13).Value = -705.9712999999997
$ws.Cells.Item(132, 14).Value = -11142.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 176.54546
$ws.Cells.Item(23, 9).Value = 80
$ws.Cells.Item(23, 10).Value = 212.75
$ws.Cells.Item(23, 11).Value = 240
$ws.Cells.Item(23, 12).Value = 638.25
$ws.Cells.Item(23, 13).Value = -5
$ws.Cells.Item(23, 14).Value = -1108.25
$ws.Cells.Item(74, 8).Value = 8600
$ws.Cells.Item(74, 10).Value = 8600
$ws.Cells.Item(74, 12).Value = 25800
$ws.Cells.Item(74, 14).Value = -27922
$ws.Cells.Item(77, 8).Value = 8600
$ws.Cells.Item(77, 10).Value = 8600
$ws.Cells.Item(77, 12).Value = 77400
$ws.Cells.Item(77, 14).Value = -88008
$ws.Cells.Item(80, 8).Value = 2500
$ws.Cells.Item(80, 10).Value = 2500
$ws.Cells.Item(80, 12).Value = 7500
$ws.Cells.Item(80, 14).Value = -9372
$ws.Cells.Item(83, 8).Value = 2500
$ws.Cells.Item(83, 10).Value = 2500
$ws.Cells.Item(83, 12).Value = 22500
$ws.Cells.Item(83, 14).Value = -31860

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5798.4
$ws.Cells.Item(70, 9).Value = 5931.467
$ws.Cells.Item(70, 10).Value = 5000
$ws.Cells.Item(70, 11).Value = 5931.467
$ws.Cells.Item(70, 12).Value = 5000
$ws.Cells.Item(70, 13).Value = -5661.467
$ws.Cells.Item(70, 14).Value = -5540
$ws.Cells.Item(73, 8).Value = 5798.4
$ws.Cells.Item(73, 9).Value = 5931.467
$ws.Cells.Item(73, 10).Value = 5000
$ws.Cells.Item(73, 11).Value = 5931.467
$ws.Cells.Item(73, 12).Value = 5000
$ws.Cells.Item(73, 13).Value = -4995.467
$ws.Cells.Item(73, 14).Value = -6872
$ws.Cells.Item(97, 8).Value = 483.05264
$ws.Cells.Item(97, 9).Value = 499.29413
$ws.Cells.Item(97, 11).Value = 499.29413
$ws.Cells.Item(97, 13).Value = -3.294129999999996
$ws.Cells.Item(113, 8).Value = 12800.5
$ws.Cells.Item(113, 9).Value = 991
$ws.Cells.Item(113, 10).Value = 14112.667
$ws.Cells.Item(113, 11).Value = 991
$ws.Cells.Item(113, 12).Value = 14112.667
$ws.Cells.Item(113, 13).Value = 1179
$ws.Cells.Item(113, 14).Value = -18452.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1208.0588
$ws.Cells.Item(61, 9).Value = 1011.8889
$ws.Cells.Item(61, 10).Value = 1428.75
$ws.Cells.Item(61, 11).Value = 1011.8889
$ws.Cells.Item(61, 12).Value = 1428.75
$ws.Cells.Item(61, 13).Value = -809.8889
$ws.Cells.Item(61, 14).Value = -1832.75
$ws.Cells.Item(113, 8).Value = 1208.0588
$ws.Cells.Item(113, 9).Value = 1011.8889
$ws.Cells.Item(113, 10).Value = 1428.75
$ws.Cells.Item(113, 11).Value = 1011.8889
$ws.Cells.Item(113, 12).Value = 1428.75
$ws.Cells.Item(113, 13).Value = 1158.1111
$ws.Cells.Item(113, 14).Value = -5768.75
$ws.Cells.Item(122, 8).Value = 4487.5
$ws.Cells.Item(122, 9).Value = 4797.5
$ws.Cells.Item(122, 10).Value = 4177.5
$ws.Cells.Item(122, 11).Value = 14392.5
$ws.Cells.Item(122, 12).Value = 12532.5
$ws.Cells.Item(122, 13).Value = -11942.5
$ws.Cells.Item(122, 14).Value = -17432.5
$ws.Cells.Item(132, 8).Value = 1505.05
$ws.Cells.Item(132, 9).Value = 1562.9468
$ws.Cells.Item(132, 10).Value = 598
$ws.Cells.Item(132, 11).Value = 4688.8404
$ws.Cells.Item(132, 12).Value = 1794
$ws.Cells.Item(132, 13).Value = -2158.8404
$ws.Cells.Item(132, 14).Value = -6854
$ws.Cells.Item(136, 8).Value = 2065.127
$ws.Cells.Item(136, 9).Value = 1648.4348
$ws.Cells.Item(136, 10).Value = 3192.647
$ws.Cells.Item(136, 11).Value = 4945.3044
$ws.Cells.Item(136, 12).Value = 9577.940999999999
$ws.Cells.Item(136, 13).Value = -2395.3044
$ws.Cells.Item(136, 14).Value = -14677.941

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 656.46155
$ws.Cells.Item(100, 9).Value = 555.5
$ws.Cells.Item(100, 10).Value = 818
$ws.Cells.Item(100, 11).Value = 1111
$ws.Cells.Item(100, 12).Value = 1636
$ws.Cells.Item(100, 13).Value = -570
$ws.Cells.Item(100, 14).Value = -2718
$ws.Cells.Item(136, 8).Value = 2404.7207
$ws.Cells.Item(136, 9).Value = 2292.491
$ws.Cells.Item(136, 10).Value = 2879.5386
$ws.Cells.Item(136, 11).Value = 6877.473
$ws.Cells.Item(136, 12).Value = 8638.6158
$ws.Cells.Item(136, 13).Value = -4327.473
$ws.Cells.Item(136, 14).Value = -13738.6158
